$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the Test Case IDs in column B for rows 20-41 sequentially
# from BANK_SYS_TC_Log_R012 to BANK_SYS_TC_Log_R033
for ($row = 20; $row -le 41; $row++) {
    $num = $row - 20 + 12
    $id = "BANK_SYS_TC_Log_R{0:D3}" -f $num
    $ws.Cells.Item($row, 2).Value = $id
}

# Update the view: scroll back to column A, and change the selection
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B9:B41").Select()
